$wb = $excel.ActiveWorkbook

# Add a new "Comments" column (column E) header to the four history sheets.
$ws2 = $wb.Worksheets.Item("Withdraw History")
$ws3 = $wb.Worksheets.Item("Deposit History")
$ws4 = $wb.Worksheets.Item("Transfer History")
$ws5 = $wb.Worksheets.Item("Absolute History")

$ws2.Range("E1").Value = "Comments"
$ws3.Range("E1").Value = "Comments"
$ws4.Range("E1").Value = "Comments"
$ws5.Range("E1").Value = "Comments"

# Restore each sheet's own selection/active-cell state.
[void]$ws3.Range("E1").Select()
[void]$ws4.Range("E1").Select()
[void]$ws5.Range("E5").Select()

# "Withdraw History" ends up the active tab/selection (matches activeTab=1).
[void]$ws2.Activate()
[void]$ws2.Range("E1").Select()
